# Update the "code list" worksheet so the customer-service related rows
# use the new "CS" category code and English wording, and the "kode input"
# column values drop their dot separators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-9: category "P" (konsumen/customer) -> "CS"
$csRows = 2..9
foreach ($r in $csRows) {
    $ws.Cells.Item($r, 2).Value = "CS"
    $ws.Cells.Item($r, 5).Value = "customer"
}

# "kode input" (column G) updates for rows 2-17
$kodeInput = @{
    2  = "ACS10"
    3  = "ACS20"
    4  = "ACS3"
    5  = "ACS4"
    6  = "ACS5"
    7  = "ACS6"
    8  = "ACS7"
    9  = "ACS8"
    10 = "APO2"
    11 = "APO3"
    12 = "AS1"
    13 = "AS2"
    14 = "AS3"
    15 = "BNF1"
    16 = "BNF2"
    17 = "BNF3"
}

foreach ($r in $kodeInput.Keys) {
    $ws.Cells.Item($r, 7).Value = $kodeInput[$r]
}

# Rows 12-14: name category "layanan" -> "service"
$serviceRows = 12..14
foreach ($r in $serviceRows) {
    $ws.Cells.Item($r, 5).Value = "service"
}
